# Add duplicate detection for contract note imports
#
# A new contract-note line (dated 2026-02-10, the CN#252611730667 buy of 20
# shares @148.71) needs to be recorded above the existing most-recent entry
# (2026-02-09, CN#252611665409, 40 shares @152.24) on the "Trading History"
# sheet, newest-first. Insert a fresh row 5, push the old row 5 down to
# row 6, then fill in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new blank row 5 - this shifts the old row 5 (and everything
# below it) down to row 6, carrying its values/formulas/styles along.
$ws.Rows("5:5").Insert()

# The inserted row inherits the header row's formatting; strip that back
# to an unformatted row, then drop the leftover (unused) cells in the
# K:AB range that Insert() also materialized so they don't linger as
# empty styled cells.
$ws.Range("A5:AB5").ClearFormats()
$ws.Range("K5:AB5").Clear()

# DATE column keeps the same custom date format used by the row below it.
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A5").Value = 46063
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 148.71
$ws.Range("F5").Value = 2995.19
$ws.Range("G5").Value = "CN#252611730667"
$ws.Range("H5").Value = 2.97
$ws.Range("I5").Value = 18.02
$ws.Range("J5").Formula = '=Index!$C$2'
